$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (52) with the July 21 raw/clean SSA data, following the
# same layout as the existing rows: A = date (stored as text, like
# "2020-07-20" etc.), B = Confirmados, C = Negativos, D = Sospechosos,
# E = Defunciones, F = Porcentaje hospitalizados.

# Column A holds date-looking strings as plain text in this sheet (e.g.
# "2020-07-20"), not real dates. Force the cell to Text format first so
# Excel doesn't auto-convert "2020-07-21" into a date serial number, then
# reset the cell style back to Normal/General so no stray formatting is
# left behind on the new cell (matching the plain, unstyled cells used by
# every other date entry in column A).
$dateCell = $ws.Cells.Item(52, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2020-07-21"
$dateCell.Style = "Normal"

$ws.Cells.Item(52, 2).Value = 356255
$ws.Cells.Item(52, 3).Value = 406151
$ws.Cells.Item(52, 4).Value = 82866
$ws.Cells.Item(52, 5).Value = 40400
$ws.Cells.Item(52, 6).Value = 28.33
